# Adds USB Host HID related translation text rows (B13:F39) to the
# "Translation" worksheet, per commit "USB Host HID 추가".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(13, 2).Value = "SingleUseId12"
$ws.Cells.Item(13, 3).Value = "Default"
$ws.Cells.Item(13, 4).Value = "Center"
$ws.Cells.Item(13, 5).Value = "USB HID"
$ws.Cells.Item(13, 6).Value = "LTR"

$ws.Cells.Item(14, 2).Value = "SingleUseId13"
$ws.Cells.Item(14, 3).Value = "SingleUseId4"
$ws.Cells.Item(14, 4).Value = "Left"
$ws.Cells.Item(14, 5).Value = "L_X : <value>"
$ws.Cells.Item(14, 6).Value = "LTR"

$ws.Cells.Item(15, 2).Value = "SingleUseId14"
$ws.Cells.Item(15, 3).Value = "SingleUseId4"
$ws.Cells.Item(15, 4).Value = "Left"
$ws.Cells.Item(15, 5).Value = "'0    "
$ws.Cells.Item(15, 6).Value = "LTR"

$ws.Cells.Item(16, 2).Value = "SingleUseId15"
$ws.Cells.Item(16, 3).Value = "SingleUseId4"
$ws.Cells.Item(16, 4).Value = "Left"
$ws.Cells.Item(16, 5).Value = "L_Y : <value>"
$ws.Cells.Item(16, 6).Value = "LTR"

$ws.Cells.Item(17, 2).Value = "SingleUseId16"
$ws.Cells.Item(17, 3).Value = "SingleUseId4"
$ws.Cells.Item(17, 4).Value = "Left"
$ws.Cells.Item(17, 5).Value = "'0    "
$ws.Cells.Item(17, 6).Value = "LTR"

$ws.Cells.Item(18, 2).Value = "SingleUseId17"
$ws.Cells.Item(18, 3).Value = "SingleUseId4"
$ws.Cells.Item(18, 4).Value = "Left"
$ws.Cells.Item(18, 5).Value = "R_X : <value>"
$ws.Cells.Item(18, 6).Value = "LTR"

$ws.Cells.Item(19, 2).Value = "SingleUseId18"
$ws.Cells.Item(19, 3).Value = "SingleUseId4"
$ws.Cells.Item(19, 4).Value = "Left"
$ws.Cells.Item(19, 5).Value = "'0    "
$ws.Cells.Item(19, 6).Value = "LTR"

$ws.Cells.Item(20, 2).Value = "SingleUseId19"
$ws.Cells.Item(20, 3).Value = "SingleUseId4"
$ws.Cells.Item(20, 4).Value = "Left"
$ws.Cells.Item(20, 5).Value = "R_Y : <value>"
$ws.Cells.Item(20, 6).Value = "LTR"

$ws.Cells.Item(21, 2).Value = "SingleUseId20"
$ws.Cells.Item(21, 3).Value = "SingleUseId4"
$ws.Cells.Item(21, 4).Value = "Left"
$ws.Cells.Item(21, 5).Value = "'0    "
$ws.Cells.Item(21, 6).Value = "LTR"

$ws.Cells.Item(22, 2).Value = "SingleUseId21"
$ws.Cells.Item(22, 3).Value = "SingleUseId4"
$ws.Cells.Item(22, 4).Value = "Left"
$ws.Cells.Item(22, 5).Value = "L_U : <value>"
$ws.Cells.Item(22, 6).Value = "LTR"

$ws.Cells.Item(23, 2).Value = "SingleUseId22"
$ws.Cells.Item(23, 3).Value = "SingleUseId4"
$ws.Cells.Item(23, 4).Value = "Left"
$ws.Cells.Item(23, 5).Value = "'0    "
$ws.Cells.Item(23, 6).Value = "LTR"

$ws.Cells.Item(24, 2).Value = "SingleUseId23"
$ws.Cells.Item(24, 3).Value = "SingleUseId4"
$ws.Cells.Item(24, 4).Value = "Left"
$ws.Cells.Item(24, 5).Value = "L_D : <value>"
$ws.Cells.Item(24, 6).Value = "LTR"

$ws.Cells.Item(25, 2).Value = "SingleUseId24"
$ws.Cells.Item(25, 3).Value = "SingleUseId4"
$ws.Cells.Item(25, 4).Value = "Left"
$ws.Cells.Item(25, 5).Value = "'0    "
$ws.Cells.Item(25, 6).Value = "LTR"

$ws.Cells.Item(26, 2).Value = "SingleUseId25"
$ws.Cells.Item(26, 3).Value = "SingleUseId4"
$ws.Cells.Item(26, 4).Value = "Left"
$ws.Cells.Item(26, 5).Value = "L_L : <value>"
$ws.Cells.Item(26, 6).Value = "LTR"

$ws.Cells.Item(27, 2).Value = "SingleUseId26"
$ws.Cells.Item(27, 3).Value = "SingleUseId4"
$ws.Cells.Item(27, 4).Value = "Left"
$ws.Cells.Item(27, 5).Value = "'0    "
$ws.Cells.Item(27, 6).Value = "LTR"

$ws.Cells.Item(28, 2).Value = "SingleUseId27"
$ws.Cells.Item(28, 3).Value = "SingleUseId4"
$ws.Cells.Item(28, 4).Value = "Left"
$ws.Cells.Item(28, 5).Value = "L_R : <value>"
$ws.Cells.Item(28, 6).Value = "LTR"

$ws.Cells.Item(29, 2).Value = "SingleUseId28"
$ws.Cells.Item(29, 3).Value = "SingleUseId4"
$ws.Cells.Item(29, 4).Value = "Left"
$ws.Cells.Item(29, 5).Value = "'0    "
$ws.Cells.Item(29, 6).Value = "LTR"

$ws.Cells.Item(30, 2).Value = "SingleUseId29"
$ws.Cells.Item(30, 3).Value = "SingleUseId4"
$ws.Cells.Item(30, 4).Value = "Left"
$ws.Cells.Item(30, 5).Value = "R_U : <value>"
$ws.Cells.Item(30, 6).Value = "LTR"

$ws.Cells.Item(31, 2).Value = "SingleUseId30"
$ws.Cells.Item(31, 3).Value = "SingleUseId4"
$ws.Cells.Item(31, 4).Value = "Left"
$ws.Cells.Item(31, 5).Value = "'0    "
$ws.Cells.Item(31, 6).Value = "LTR"

$ws.Cells.Item(32, 2).Value = "SingleUseId31"
$ws.Cells.Item(32, 3).Value = "SingleUseId4"
$ws.Cells.Item(32, 4).Value = "Left"
$ws.Cells.Item(32, 5).Value = "R_D : <value>"
$ws.Cells.Item(32, 6).Value = "LTR"

$ws.Cells.Item(33, 2).Value = "SingleUseId32"
$ws.Cells.Item(33, 3).Value = "SingleUseId4"
$ws.Cells.Item(33, 4).Value = "Left"
$ws.Cells.Item(33, 5).Value = "'0    "
$ws.Cells.Item(33, 6).Value = "LTR"

$ws.Cells.Item(34, 2).Value = "SingleUseId33"
$ws.Cells.Item(34, 3).Value = "SingleUseId4"
$ws.Cells.Item(34, 4).Value = "Left"
$ws.Cells.Item(34, 5).Value = "R_L : <value>"
$ws.Cells.Item(34, 6).Value = "LTR"

$ws.Cells.Item(35, 2).Value = "SingleUseId34"
$ws.Cells.Item(35, 3).Value = "SingleUseId4"
$ws.Cells.Item(35, 4).Value = "Left"
$ws.Cells.Item(35, 5).Value = "'0    "
$ws.Cells.Item(35, 6).Value = "LTR"

$ws.Cells.Item(36, 2).Value = "SingleUseId35"
$ws.Cells.Item(36, 3).Value = "SingleUseId4"
$ws.Cells.Item(36, 4).Value = "Left"
$ws.Cells.Item(36, 5).Value = "R_R : <value>"
$ws.Cells.Item(36, 6).Value = "LTR"

$ws.Cells.Item(37, 2).Value = "SingleUseId36"
$ws.Cells.Item(37, 3).Value = "SingleUseId4"
$ws.Cells.Item(37, 4).Value = "Left"
$ws.Cells.Item(37, 5).Value = "'0    "
$ws.Cells.Item(37, 6).Value = "LTR"

$ws.Cells.Item(38, 2).Value = "SingleUseId37"
$ws.Cells.Item(38, 3).Value = "SingleUseId4"
$ws.Cells.Item(38, 4).Value = "Left"
$ws.Cells.Item(38, 5).Value = "USB : <value>"
$ws.Cells.Item(38, 6).Value = "LTR"

$ws.Cells.Item(39, 2).Value = "SingleUseId38"
$ws.Cells.Item(39, 3).Value = "SingleUseId4"
$ws.Cells.Item(39, 4).Value = "Left"
$ws.Cells.Item(39, 5).Value = "Disconnected "
$ws.Cells.Item(39, 6).Value = "LTR"
